$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1274.375
$ws.Range("I2").Value = 32.666668
$ws.Range("J2").Value = 4999.5
$ws.Range("K2").Value = 32.666668
$ws.Range("L2").Value = 4999.5
$ws.Range("M2").Value = 80.333332
$ws.Range("N2").Value = -5225.5

$ws.Range("H28").Value = 676.5
$ws.Range("I28").Value = 511.8
$ws.Range("J28").Value = 1500
$ws.Range("K28").Value = 511.8
$ws.Range("L28").Value = 1500
$ws.Range("M28").Value = -26.80000000000001
$ws.Range("N28").Value = -2470

$ws.Range("H41").Value = 847.9
$ws.Range("I41").Value = 507.25
$ws.Range("J41").Value = 1075
$ws.Range("K41").Value = 507.25
$ws.Range("L41").Value = 1075
$ws.Range("M41").Value = -67.25
$ws.Range("N41").Value = -1955

$ws.Range("H43").Value = 1991.75
$ws.Range("I43").Value = 1988.3334
$ws.Range("K43").Value = 1988.3334
$ws.Range("M43").Value = -1919.3334

$ws.Range("H53").Value = 516.7
$ws.Range("I53").Value = 710.4286
$ws.Range("J53").Value = 64.666664
$ws.Range("K53").Value = 710.4286
$ws.Range("L53").Value = 64.666664
$ws.Range("M53").Value = -73.42859999999996
$ws.Range("N53").Value = -1338.666664

$ws.Range("H62").Value = 5475
$ws.Range("I62").Value = 4791.375
$ws.Range("J62").Value = 6256.2856
$ws.Range("K62").Value = 4791.375
$ws.Range("L62").Value = 6256.2856
$ws.Range("M62").Value = -4167.375
$ws.Range("N62").Value = -7504.2856

$ws.Range("H65").Value = 5475
$ws.Range("I65").Value = 4791.375
$ws.Range("J65").Value = 6256.2856
$ws.Range("K65").Value = 23956.875
$ws.Range("L65").Value = 31281.428
$ws.Range("M65").Value = -20836.875
$ws.Range("N65").Value = -37521.428

$ws.Range("H76").Value = 5212.1665
$ws.Range("I76").Value = 4759.6665
$ws.Range("K76").Value = 4759.6665
$ws.Range("M76").Value = -4444.6665

$ws.Range("H79").Value = 5212.1665
$ws.Range("I79").Value = 4759.6665
$ws.Range("K79").Value = 4759.6665
$ws.Range("M79").Value = -3667.6665

$ws.Range("H86").Value = 16749.834
$ws.Range("I86").Value = 21249.75
$ws.Range("J86").Value = 7750
$ws.Range("K86").Value = 21249.75
$ws.Range("L86").Value = 7750
$ws.Range("M86").Value = -20126.75
$ws.Range("N86").Value = -9996

$ws.Range("H89").Value = 16749.834
$ws.Range("I89").Value = 21249.75
$ws.Range("J89").Value = 7750
$ws.Range("K89").Value = 106248.75
$ws.Range("L89").Value = 38750
$ws.Range("M89").Value = -100632.75
$ws.Range("N89").Value = -49982

$ws.Range("H92").Value = 598.3333
$ws.Range("I92").Value = 1295
$ws.Range("K92").Value = 1295
$ws.Range("M92").Value = -47

$ws.Range("H98").Value = 3522.5833
$ws.Range("I98").Value = 3553
$ws.Range("K98").Value = 3553
$ws.Range("M98").Value = -2055

$ws.Range("H106").Value = 4800
$ws.Range("I106").Value = 4800
$ws.Range("K106").Value = 4800
$ws.Range("M106").Value = -4169

$ws.Range("H107").Value = 3441.5
$ws.Range("J107").Value = 10000
$ws.Range("L107").Value = 10000
$ws.Range("N107").Value = -13840

$ws.Range("H122").Value = 3522.5833
$ws.Range("I122").Value = 3553
$ws.Range("K122").Value = 10659
$ws.Range("M122").Value = -8209

$ws.Range("H131").Value = 1342
$ws.Range("I131").Value = 998.8
$ws.Range("J131").Value = 2200
$ws.Range("K131").Value = 2996.4
$ws.Range("L131").Value = 6600
$ws.Range("M131").Value = 2043.6
$ws.Range("N131").Value = -16680

$ws.Range("H137").Value = 1410.3334
$ws.Range("I137").Value = 1167.909
$ws.Range("J137").Value = 2477
$ws.Range("K137").Value = 3503.727
$ws.Range("L137").Value = 7431
$ws.Range("M137").Value = -953.7270000000003
$ws.Range("N137").Value = -12531

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 275
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 275
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 275
$ws.Range("M3").ClearContents()
$ws.Range("N3").Value = -505

$ws.Range("H5").Value = 272.2857
$ws.Range("I5").Value = 272.2857
$ws.Range("K5").Value = 272.2857
$ws.Range("M5").Value = -160.2857

$ws.Range("H32").Value = 5539.675
$ws.Range("I32").Value = 4357.5527
$ws.Range("K32").Value = 4357.5527
$ws.Range("M32").Value = -4070.5527

$ws.Range("H36").Value = 11250
$ws.Range("I36").Value = 3250
$ws.Range("J36").Value = 19250
$ws.Range("K36").Value = 3250
$ws.Range("L36").Value = 19250
$ws.Range("M36").Value = -2904
$ws.Range("N36").Value = -19942

$ws.Range("H55").Value = 29000
$ws.Range("J55").Value = 29000
$ws.Range("L55").Value = 29000
$ws.Range("N55").Value = -29630

$ws.Range("H61").Value = 2141.7144
$ws.Range("I61").Value = 1348.4
$ws.Range("K61").Value = 1348.4
$ws.Range("M61").Value = -1136.4

$ws.Range("H136").Value = 2141.7144
$ws.Range("I136").Value = 1348.4
$ws.Range("K136").Value = 4045.2
$ws.Range("M136").Value = -1495.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 272.2857
$ws.Range("I4").Value = 272.2857
$ws.Range("K4").Value = 272.2857
$ws.Range("M4").Value = -157.2857

$ws.Range("H22").Value = 256
$ws.Range("I22").Value = 260
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 260
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = -87
$ws.Range("N22").Value = -596

$ws.Range("H86").Value = 4300
$ws.Range("I86").Value = 4300
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 4300
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -3177
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 4300
$ws.Range("I89").Value = 4300
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 21500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -15884
$ws.Range("N89").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 726.6667
$ws.Range("J22").Value = 726.6667
$ws.Range("L22").Value = 726.6667
$ws.Range("N22").Value = -1426.6667

$ws.Range("H94").Value = 2578.4
$ws.Range("I94").Value = 2486
$ws.Range("K94").Value = 2486
$ws.Range("M94").Value = -2035

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 350
$ws.Range("I87").Value = 350
$ws.Range("K87").Value = 1050
$ws.Range("M87").Value = 198

$ws.Range("H90").Value = 350
$ws.Range("I90").Value = 350
$ws.Range("K90").Value = 3150
$ws.Range("M90").Value = 3090

$ws.Range("H132").Value = 3377.5
$ws.Range("I132").Value = 3255
$ws.Range("K132").Value = 29295
$ws.Range("M132").Value = -26765

$ws.Range("H139").Value = 2499.5
$ws.Range("I139").Value = 2499.5
$ws.Range("K139").Value = 7498.5
$ws.Range("M139").Value = -2358.5

$ws.Range("H140").Value = 1569.3334
$ws.Range("I140").Value = 963.2
$ws.Range("J140").Value = 4600
$ws.Range("K140").Value = 2889.6
$ws.Range("L140").Value = 13800
$ws.Range("M140").Value = 2290.4
$ws.Range("N140").Value = -24160

$ws.Range("H141").Value = 1999
$ws.Range("I141").Value = 1999
$ws.Range("K141").Value = 5997
$ws.Range("M141").Value = -817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 20000
$ws.Range("J14").Value = 20000
$ws.Range("L14").Value = 20000
$ws.Range("N14").Value = -20336

$ws.Range("H49").Value = 31000
$ws.Range("J49").Value = 31000
$ws.Range("L49").Value = 31000
$ws.Range("N49").Value = -31368

$ws.Range("H126").Value = 2251
$ws.Range("I126").Value = 1668
$ws.Range("K126").Value = 5004
$ws.Range("M126").Value = -2534

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 4500
$ws.Range("I5").Value = 4500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 4500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -4387
$ws.Range("N5").ClearContents()

$ws.Range("H22").Value = 1043.75
$ws.Range("J22").Value = 1475
$ws.Range("L22").Value = 1475
$ws.Range("N22").Value = -2065

$ws.Range("H27").Value = 1043.75
$ws.Range("J27").Value = 1475
$ws.Range("L27").Value = 1475
$ws.Range("N27").Value = -1689

$ws.Range("H82").Value = 1632.2222
$ws.Range("I82").Value = 1472.5
$ws.Range("J82").Value = 1760
$ws.Range("K82").Value = 1472.5
$ws.Range("L82").Value = 1760
$ws.Range("M82").Value = -1111.5
$ws.Range("N82").Value = -2482

$ws.Range("H85").Value = 1632.2222
$ws.Range("I85").Value = 1472.5
$ws.Range("J85").Value = 1760
$ws.Range("K85").Value = 1472.5
$ws.Range("L85").Value = 1760
$ws.Range("M85").Value = -224.5
$ws.Range("N85").Value = -4256

$ws.Range("H95").Value = 19999
$ws.Range("J95").Value = 19999
$ws.Range("L95").Value = 19999
$ws.Range("N95").Value = -25491

$ws.Range("H122").Value = 3797.7646
$ws.Range("I122").Value = 3310.3635
$ws.Range("J122").Value = 4691.3335
$ws.Range("K122").Value = 9931.0905
$ws.Range("L122").Value = 14074.0005
$ws.Range("M122").Value = -7481.0905
$ws.Range("N122").Value = -18974.0005

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 9999999
$ws.Range("I15").Value = 9999999
$ws.Range("K15").Value = 9999999
$ws.Range("M15").Value = -9999711

$ws.Range("H31").Value = 519950
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 519950
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 519950
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -520646

$ws.Range("H41").Value = 21259.2
$ws.Range("I41").Value = 21565.666
$ws.Range("J41").Value = 20799.5
$ws.Range("K41").Value = 21565.666
$ws.Range("L41").Value = 20799.5
$ws.Range("M41").Value = -21175.666
$ws.Range("N41").Value = -21579.5

$ws.Range("H132").Value = 1625.75
$ws.Range("I132").Value = 1601
$ws.Range("K132").Value = 4803
$ws.Range("M132").Value = -2273
